$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 (was PYPL 79.45 / 15) to MSFT data (416.61 / 3)
$ws.Range("A5").Value = "MSFT"
$ws.Range("B5").Value = 416.61
$ws.Range("C5").Value = 3

# Remove the now-duplicated rows 6 (NVDA) and 7 (old MSFT) entirely,
# shifting everything up so the table ends at row 5
$ws.Range("A6:C7").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Move the active selection to reflect the new last empty row under the table
$ws.Range("C4").Select()
